$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row 4, mirroring row 3's layout/values but with an updated date/time in column A.
$ws.Range("A4").Value = 42602.01458333333
$ws.Range("A4").NumberFormat = "m/d/yy h:mm"

$ws.Range("B4").Value = "Clean Energy Fuels Corp."
$ws.Range("C4").Value = "CLNE"
$ws.Range("D4").Value = 4.1500000000000004
$ws.Range("E4").Value = 4.3499999999999996
$ws.Range("F4").Value = 6.88
$ws.Range("G4").Value = 4.07
